$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-34 with corrected data (dates fixed from April to March, plus other corrections)
$data = @(
    @(44650, 'Oriol', 'Máximo Agustí Galiano', '96746138L', 'terapia', 53),
    @(44638, 'Elena', 'Lola Conesa Agustí', '92116835H', 'terapia', 64),
    @(44639, 'Elena', 'Lola Conesa Agustí', '92116835H', 'terapia', 68),
    @(44648, 'Elena', 'Lola Conesa Agustí', '92116835H', 'terapia', 50),
    @(44626, 'Elena', 'Lola Conesa Agustí', '92116835H', 'terapia', 75),
    @(44633, 'Elena', 'Lola Conesa Agustí', '92116835H', 'terapia', 71),
    @(44636, 'Elena', 'Lola Conesa Agustí', '92116835H', 'terapia', 76),
    @(44649, 'Elena', 'Lola Conesa Agustí', '92116835H', 'terapia', 60),
    @(44641, 'Oriol', 'Serafina Hernandez Blanch', '49144449G', 'terapia', 58),
    @(44648, 'Oriol', 'Serafina Hernandez Blanch', '49144449G', 'terapia', 56),
    @(44630, 'Oriol', 'Serafina Hernandez Blanch', '49144449G', 'terapia', 79),
    @(44624, 'Oriol', 'Serafina Hernandez Blanch', '49144449G', 'terapia', 66),
    @(44640, 'Oriol', 'Serafina Hernandez Blanch', '49144449G', 'terapia', 73),
    @(44622, 'Oriol', 'Serafina Hernandez Blanch', '49144449G', 'terapia', 56),
    @(44649, 'Oriol', 'Serafina Hernandez Blanch', '49144449G', 'terapia', 62),
    @(44650, 'Oriol', 'Serafina Hernandez Blanch', '49144449G', 'terapia', 74),
    @(44637, 'Oriol', 'Serafina Hernandez Blanch', '49144449G', 'terapia', 72),
    @(44643, 'Oriol', 'Serafina Hernandez Blanch', '49144449G', 'terapia', 58),
    @(44634, 'Oriol', 'Osvaldo Nogués Palacios', '33841962S', 'terapia', 80),
    @(44627, 'Psicologo3', 'Cristian Carbajo Hernandez', '25820527Z', 'terapia', 54),
    @(44647, 'Psicologo3', 'Cristian Carbajo Hernandez', '25820527Z', 'terapia', 69),
    @(44640, 'Psicologo3', 'Cristian Carbajo Hernandez', '25820527Z', 'terapia', 52),
    @(44628, 'Psicologo3', 'Cristian Carbajo Hernandez', '25820527Z', 'terapia', 64),
    @(44625, 'Psicologo3', 'Cristian Carbajo Hernandez', '25820527Z', 'terapia', 63),
    @(44630, 'Psicologo3', 'Cristian Carbajo Hernandez', '25820527Z', 'terapia', 80),
    @(44640, 'Psicologo3', 'Cristian Carbajo Hernandez', '25820527Z', 'terapia', 52),
    @(44639, 'Psicologo3', 'Cristian Carbajo Hernandez', '25820527Z', 'terapia', 69),
    @(44637, 'Oriol', 'Ciriaco de Acosta', '74962378L', 'terapia', 56),
    @(44646, 'Oriol', 'Ciriaco de Acosta', '74962378L', 'terapia', 68),
    @(44643, 'Oriol', 'Ciriaco de Acosta', '74962378L', 'terapia', 54),
    @(44644, 'Oriol', 'Ciriaco de Acosta', '74962378L', 'terapia', 69),
    @(44651, 'Oriol', 'Ciriaco de Acosta', '74962378L', 'terapia', 54),
    @(44633, 'Oriol', 'Ciriaco de Acosta', '74962378L', 'terapia', 64)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Remove now-obsolete rows 35-47 (table shrank from 47 to 34 data-bearing rows)
$ws.Range("A35:F47").EntireRow.Delete()

Write-Host "Done. New used range:" $ws.UsedRange.Address